$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 1.17
$ws.Range("N2").Value = 5
$ws.Range("O2").Value = 1.8
$ws.Range("P2").Value = 1.91
$ws.Range("AB2").Value = 1.4
$ws.Range("AA3").Value = 1.57
$ws.Range("G5").Value = 2.62
$ws.Range("H5").Value = 3.25
$ws.Range("I5").Value = 2.5
$ws.Range("J5").Value = 3.2
$ws.Range("K5").Value = 2.07
$ws.Range("L5").Value = 3.05
$ws.Range("O5").Value = 1.26
$ws.Range("P5").Value = 3.2
$ws.Range("S5").Value = 1.75
$ws.Range("T5").Value = 1.85
$ws.Range("W5").Value = 2.77
$ws.Range("X5").Value = 1.34
$ws.Range("AA5").Value = 1.6
$ws.Range("AB5").Value = 2.05
$ws.Range("AC5").Value = 9
$ws.Range("AD5").Value = 14
$ws.Range("AG5").Value = 21
$ws.Range("AH5").Value = 28
$ws.Range("AI5").Value = 10.5
$ws.Range("AJ5").Value = 6.4
$ws.Range("AK5").Value = 12.5
$ws.Range("AL5").Value = 50
$ws.Range("AM5").Value = 9.5
$ws.Range("AN5").Value = 13.5
$ws.Range("AQ5").Value = 19.5
$ws.Range("AR5").Value = 26
$ws.Range("AS5").Value = 350
$ws.Range("G7").Value = 2.35
$ws.Range("M7").Value = 1.04
$ws.Range("N7").Value = 9
$ws.Range("S7").Value = 2.08
$ws.Range("T7").Value = 1.73
$ws.Range("W7").Value = 3.75
$ws.Range("X7").Value = 1.25
$ws.Range("AF7").Value = 23
$ws.Range("AH7").Value = 34
$ws.Range("AO7").Value = 11
$ws.Range("G8").Value = 1.75
$ws.Range("H8").Value = 3.7
$ws.Range("I8").Value = 3.9
$ws.Range("J8").Value = 2.38
$ws.Range("K8").Value = 2.2
$ws.Range("M8").Value = 1.04
$ws.Range("N8").Value = 9
$ws.Range("O8").Value = 1.25
$ws.Range("P8").Value = 3.75
$ws.Range("S8").Value = 1.8
$ws.Range("T8").Value = 2
$ws.Range("W8").Value = 3
$ws.Range("X8").Value = 1.36
$ws.Range("Y8").Value = 1.36
$ws.Range("Z8").Value = 3
$ws.Range("AA8").Value = 1.8
$ws.Range("AB8").Value = 1.91
$ws.Range("AC8").Value = 7.5
$ws.Range("AH8").Value = 26
$ws.Range("AI8").Value = 11
$ws.Range("AJ8").Value = 7.5
$ws.Range("AM8").Value = 12
$ws.Range("AS8").Value = 500
